$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coordinates")

$ws.Range("D2").Value = 32.094662

$ws.Range("C3").Value = -31.638883
$ws.Range("D3").Value = 30.870622

$ws.Range("C4").Value = -33.318683
$ws.Range("D4").Value = 29.505758

$ws.Range("C5").Value = -34.895537
$ws.Range("D5").Value = 27.046444

$ws.Range("C6").Value = -35.940604
$ws.Range("D6").Value = 23.644859

$ws.Range("C7").Value = -36.169229
$ws.Range("D7").Value = 19.822039

$ws.Range("C8").Value = -34.685275
$ws.Range("D8").Value = 16.876437

$ws.Range("C9").Value = -32.553781
$ws.Range("D9").Value = 14.83723

$ws.Range("C10").Value = -30.469595
$ws.Range("D10").Value = 12.853821

$ws.Range("C11").Value = -27.293171
$ws.Range("D11").Value = 11.269671

$ws.Range("C12").Value = -24.474346
$ws.Range("D12").Value = 10.365769

$ws.Range("C13").Value = -21.280331
$ws.Range("D13").Value = 9.349475

$ws.Range("C14").Value = -17.493437
$ws.Range("D14").Value = 8.281086

$ws.Range("C15").Value = -14.646522
$ws.Range("D15").Value = 8.102
